$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataEntry")

# Fix trailing-slash typos in the two IEEE citation URLs (B13, then B12)
$ws.Range("B13").Value = "Lumpkins, W. (n.d.). The MobiAria Wireless Bluetooth Speaker. Retrieved September 18, 2017, from http://ieeexplore.ieee.org/document/6685931`nBodson, D. (n.d.). Digital Audio Around the World. Retrieved from http://ieeexplore.ieee.org/stamp/stamp.jsp?tp=&arnumber=5641649`nPauli, M. (2017, May 5). Miniaturized Millimeter-Wave Radar Sensor for High-Accuracy Applications. Retrieved from http://ieeexplore.ieee.org/stamp/stamp.jsp?tp=&arnumber=7885501               "

$ws.Range("B12").Value = "Bluetooth in wireless communication. (n.d.). Retrieved September 18, 2017, from http://ieeexplore.ieee.org/document/1007414`nSparkFun FM Tuner Evaluation Board - Si4703. (n.d.). Retrieved from https://www.sparkfun.com/products/12938`n"
